$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.47"
$ws.Range("D3").Value = "'21.93"
$ws.Range("D4").Value = "'5.387"
$ws.Range("D5").Value = "'0.05858"
$ws.Range("D7").Value = "'6.366"
$ws.Range("D8").Value = "'0.8136"
$ws.Range("D9").Value = "'1.018"
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = "'0.1422"
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = "'0.03921"
$ws.Range("E11").Value = '10LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = "'0.07420"
$ws.Range("E12").Value = '11MandalaExchangeTokenMDX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = "'0.03041"
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'MCDex'
$ws.Range("C14").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D14").Value = "'4.166"
$ws.Range("E14").Value = '13MCDexMCB'
$ws.Range("B15").Value = 'BitMartToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D15").Value = "'0.09402"
$ws.Range("E15").Value = '14BitMartTokenBMX'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").Value = "'0.001602"
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").Value = "'0.04833"
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'One'
$ws.Range("C18").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D18").Value = "'0.0005891"
$ws.Range("E18").Value = '17OneONEWorstin24h'
$ws.Range("D19").Value = "'0.005799"
$ws.Range("D20").Value = "'0.004081"
$ws.Range("D21").Value = "'0.0009933"
$ws.Range("D23").Value = "'3.745"
$ws.Range("D24").Value = "'2.229"
$ws.Range("D27").Value = "'0.0002483"
$ws.Range("E27").Value = '26UpBotsUBXT'
$ws.Range("D40").Value = "'0.03874"
$ws.Range("D42").Value = "'0.1074"
$ws.Range("D43").Value = "'0.002600"
$ws.Range("E43").Value = '42CEJICEJI'
$ws.Range("D44").Value = "'0.005156"
$ws.Range("D45").Value = "'0.00005631"
$ws.Range("D47").Value = "'1.290"
$ws.Range("E47").Value = '46CoinbaseStockTokenCOINBestin24h'
$ws.Range("D48").Value = "'0.1424"
